$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric-valued cells (population/case/test counts, etc.) ---
$ws.Range("C2").Value = 864883.0
$ws.Range("F2").Value = 2565701.0
$ws.Range("I2").Value = 11287.0
$ws.Range("Y2").Value = 10375.0
$ws.Range("C3").Value = 249576.0
$ws.Range("F3").Value = 613517.0
$ws.Range("G3").Value = 33244.0
$ws.Range("I3").Value = 3316.0
$ws.Range("Y3").Value = 3106.0
$ws.Range("C4").Value = 73703.0
$ws.Range("D4").Value = 30152.0
$ws.Range("F4").Value = 194681.0
$ws.Range("G4").Value = 38448.0
$ws.Range("C5").Value = 24062.0
$ws.Range("D5").Value = 28337.0
$ws.Range("C6").Value = 105072.0
$ws.Range("D6").Value = 26713.0
$ws.Range("F6").Value = 253031.0
$ws.Range("G6").Value = 31056.0
$ws.Range("C7").Value = 46739.0
$ws.Range("D7").Value = 27781.0
$ws.Range("F7").Value = 104077.0
$ws.Range("G7").Value = 29864.0
$ws.Range("I7").Value = 758.0
$ws.Range("Y7").Value = 718.0
$ws.Range("C8").Value = 141783.0
$ws.Range("D8").Value = 20362.0
$ws.Range("F8").Value = 442659.0
$ws.Range("I8").Value = 1655.0
$ws.Range("Y8").Value = 1485.0
$ws.Range("C9").Value = 89622.0
$ws.Range("F9").Value = 327188.0
$ws.Range("C10").Value = 42954.0
$ws.Range("I11").Value = 16.0
$ws.Range("Y11").Value = 15.0
$ws.Range("C12").Value = 121685.0
$ws.Range("F12").Value = 471276.0
$ws.Range("I12").Value = 2925.0
$ws.Range("Y12").Value = 2681.0
$ws.Range("C13").Value = 54407.0
$ws.Range("F13").Value = 143244.0
$ws.Range("G13").Value = 20638.0
$ws.Range("I13").Value = 827.0
$ws.Range("Y13").Value = 770.0
$ws.Range("F14").Value = 137885.0
$ws.Range("G14").Value = 47388.0
$ws.Range("I14").Value = 773.0
$ws.Range("Y14").Value = 698.0
$ws.Range("F15").Value = 85832.0
$ws.Range("C16").Value = 24346.0
$ws.Range("D16").Value = 18176.0
$ws.Range("I16").Value = 663.0
$ws.Range("Y16").Value = 596.0
$ws.Range("C17").Value = 77371.0
$ws.Range("F17").Value = 220273.0
$ws.Range("G17").Value = 26676.0
$ws.Range("C18").Value = 40498.0
$ws.Range("D18").Value = 20149.0
$ws.Range("F21").Value = 34322.0
$ws.Range("G21").Value = 21166.0
$ws.Range("F23").Value = 34887.0
$ws.Range("C24").Value = 213686.0
$ws.Range("I24").Value = 2042.0
$ws.Range("Y24").Value = 1864.0
$ws.Range("C26").Value = 3762.0
$ws.Range("D26").Value = 14089.0
$ws.Range("F26").Value = 9589.0
$ws.Range("G26").Value = 17337.0
$ws.Range("C27").Value = 3491.0
$ws.Range("D27").Value = 17702.0
$ws.Range("F27").Value = 12546.0
$ws.Range("G27").Value = 30712.0
$ws.Range("C28").Value = 39012.0
$ws.Range("F28").Value = 104049.0
$ws.Range("I28").Value = 400.0
$ws.Range("Y28").Value = 368.0
$ws.Range("I29").Value = 277.0
$ws.Range("Y29").Value = 255.0
$ws.Range("C30").Value = 22445.0
$ws.Range("D30").Value = 16434.0
$ws.Range("F30").Value = 66258.0
$ws.Range("C31").Value = 137878.0
$ws.Range("D31").Value = 18385.0
$ws.Range("F31").Value = 411601.0
$ws.Range("G31").Value = 26496.0
$ws.Range("I31").Value = 1277.0
$ws.Range("Y31").Value = 1164.0
$ws.Range("C32").Value = 60782.0
$ws.Range("D32").Value = 22847.0
$ws.Range("C33").Value = 22482.0
$ws.Range("D33").Value = 23274.0
$ws.Range("C34").Value = 38300.0

# --- Update text-valued cells (ranges / CI strings stored as text) ---
# Force text number-format so these remain strings (matches source formatting),
# then clear the format afterwards so no stray style is left on the cell.
$textCells = @("E2", "R2", "E4", "E5", "E6", "E7", "Z7", "E8", "Z8", "E9", "E10", "E11", "R11", "Z11", "E12", "Z12", "AD12", "E13", "E14", "H14", "E15", "E16", "R16", "E19", "E20", "E21", "E23", "E25", "E26", "H26", "E27", "H27", "E28", "R28", "Z28", "E29", "Z29", "E31", "Z31", "E32", "E33", "E34")
foreach ($cellRef in $textCells) { $ws.Range($cellRef).NumberFormat = "@" }

$ws.Range("E2").Value = '1.03 (0.95-1.11)'
$ws.Range("R2").Value = '6.2 (5.8-6.7)'
$ws.Range("E4").Value = '1.01 (0.91-1.11)'
$ws.Range("E5").Value = '1.09 (0.98-1.20)'
$ws.Range("E6").Value = '1.02 (0.95-1.09)'
$ws.Range("E7").Value = '0.95 (0.85-1.05)'
$ws.Range("Z7").Value = '94.7 (92.9-96.2)'
$ws.Range("E8").Value = '1.03 (0.92-1.13)'
$ws.Range("Z8").Value = '89.7 (88.2-91.1)'
$ws.Range("E9").Value = '1.11 (1.03-1.20)'
$ws.Range("E10").Value = '1.02 (0.91-1.12)'
$ws.Range("E11").Value = '0.96 (0.83-1.09)'
$ws.Range("R11").Value = '6.2 (0.2-30.2)'
$ws.Range("Z11").Value = '93.8 (69.8-99.8)'
$ws.Range("E12").Value = '1.13 (1.03-1.22)'
$ws.Range("Z12").Value = '91.7 (90.6-92.6)'
$ws.Range("AD12").Value = '0.1 (0.0-0.3)'
$ws.Range("E13").Value = '1.11 (1.04-1.19)'
$ws.Range("E14").Value = '1.13 (1.02-1.23)'
$ws.Range("H14").Value = '19.66'
$ws.Range("E15").Value = '1.14 (1.02-1.25)'
$ws.Range("E16").Value = '1.13 (1.03-1.23)'
$ws.Range("R16").Value = '9.0 (7.0-11.5)'
$ws.Range("E19").Value = '1.08 (0.93-1.24)'
$ws.Range("E20").Value = '1.10 (0.93-1.27)'
$ws.Range("E21").Value = '1.11 (0.99-1.23)'
$ws.Range("E23").Value = '1.11 (1.00-1.22)'
$ws.Range("E25").Value = '1.08 (0.88-1.28)'
$ws.Range("E26").Value = '1.13 (0.97-1.29)'
$ws.Range("H26").Value = '43.33'
$ws.Range("E27").Value = '1.14 (0.98-1.32)'
$ws.Range("H27").Value = '31.95'
$ws.Range("E28").Value = '1.11 (1.00-1.22)'
$ws.Range("R28").Value = '6.2 (4.1-9.1)'
$ws.Range("Z28").Value = '92.0 (88.9-94.5)'
$ws.Range("E29").Value = '1.05 (0.93-1.18)'
$ws.Range("Z29").Value = '92.1 (88.2-95.0)'
$ws.Range("E31").Value = '1.06 (0.96-1.16)'
$ws.Range("Z31").Value = '91.2 (89.5-92.7)'
$ws.Range("E32").Value = '0.94 (0.84-1.03)'
$ws.Range("E33").Value = '0.94 (0.84-1.04)'
$ws.Range("E34").Value = '0.92 (0.82-1.01)'

foreach ($cellRef in $textCells) { $ws.Range($cellRef).ClearFormats() }
